# Update cryptos worksheet: Price (D) and Volume(1h) (E) columns
# Commit: Updated cryptos list on Sat Jun  1 03:52:55 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '67.737.67'
Set-TextValue $ws.Range('E2') '  -1.26%  '
Set-TextValue $ws.Range('D3') '3.779.94'
Set-TextValue $ws.Range('E3') '  +0.47%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '595.32'
Set-TextValue $ws.Range('E5') '  +0.23%  '
Set-TextValue $ws.Range('D6') '167.23'
Set-TextValue $ws.Range('E6') '  +0.01%  '
Set-TextValue $ws.Range('D7') '3.776.12'
Set-TextValue $ws.Range('E7') '  +0.45%  '
Set-TextValue $ws.Range('E8') '  +0.05%  '
Set-TextValue $ws.Range('D10') '0.159'
Set-TextValue $ws.Range('E10') '  -0.50%  '
Set-TextValue $ws.Range('D11') '6.29'
Set-TextValue $ws.Range('E11') '  -2.02%  '
Set-TextValue $ws.Range('E12') '  -0.43%  '
Set-TextValue $ws.Range('D13') '0.0000253'
Set-TextValue $ws.Range('D14') '36.01'
Set-TextValue $ws.Range('E14') '  -0.26%  '
Set-TextValue $ws.Range('D15') '4.414.86'
Set-TextValue $ws.Range('E15') '  +0.42%  '
Set-TextValue $ws.Range('D16') '3.780.69'
Set-TextValue $ws.Range('E16') '  +0.29%  '
Set-TextValue $ws.Range('D17') '67.701.28'
Set-TextValue $ws.Range('E17') '  -1.24%  '
Set-TextValue $ws.Range('D18') '18.34'
Set-TextValue $ws.Range('E18') '  +2.24%  '
Set-TextValue $ws.Range('E19') '  +0.02%  '
Set-TextValue $ws.Range('E20') '  -0.95%  '
Set-TextValue $ws.Range('D21') '10.01'
Set-TextValue $ws.Range('E21') '  -6.75%  '
Set-TextValue $ws.Range('D22') '457.28'
Set-TextValue $ws.Range('E22') '  -1.77%  '
Set-TextValue $ws.Range('D23') '0.694'
Set-TextValue $ws.Range('E23') '  -0.40%  '
Set-TextValue $ws.Range('E24') '  +3.53%  '
Set-TextValue $ws.Range('D25') '83.19'
Set-TextValue $ws.Range('E25') '  -1.21%  '
Set-TextValue $ws.Range('D26') '11.95'
Set-TextValue $ws.Range('E26') '  -0.15%  '
Set-TextValue $ws.Range('E27') '  -2.56%  '
Set-TextValue $ws.Range('D28') '10.03'
Set-TextValue $ws.Range('E28') '  -0.31%  '
Set-TextValue $ws.Range('E29') '  +0.12%  '
Set-TextValue $ws.Range('E30') '  -0.46%  '
Set-TextValue $ws.Range('E31') '  +3.06%  '
Set-TextValue $ws.Range('E32') '  -1.56%  '
Set-TextValue $ws.Range('D33') '29.69'
Set-TextValue $ws.Range('E33') '  -1.08%  '
Set-TextValue $ws.Range('D34') '9.12'
Set-TextValue $ws.Range('E34') '  -0.95%  '
Set-TextValue $ws.Range('E35') '  +0.17%  '
Set-TextValue $ws.Range('D36') '3.731.74'
Set-TextValue $ws.Range('E36') '  +0.45%  '
Set-TextValue $ws.Range('E37') '  -0.55%  '
Set-TextValue $ws.Range('D38') '3.33'
Set-TextValue $ws.Range('E38') '  -2.27%  '
Set-TextValue $ws.Range('E39') '  -0.54%  '
Set-TextValue $ws.Range('D40') '0.995'
Set-TextValue $ws.Range('E40') '  -0.34%  '
Set-TextValue $ws.Range('D41') '5.76'
Set-TextValue $ws.Range('E41') '  -0.69%  '
Set-TextValue $ws.Range('D42') '1.00'
Set-TextValue $ws.Range('E42') '  -0.02%  '
Set-TextValue $ws.Range('D44') '45.38'
Set-TextValue $ws.Range('E44') '  +3.37%  '
Set-TextValue $ws.Range('D45') '48.22'
Set-TextValue $ws.Range('E45') '  +3.00%  '
Set-TextValue $ws.Range('D46') '0.298'
Set-TextValue $ws.Range('E46') '  -1.40%  '
Set-TextValue $ws.Range('D47') '149.71'
Set-TextValue $ws.Range('E47') '  +3.12%  '
Set-TextValue $ws.Range('D48') '8.31'
Set-TextValue $ws.Range('E48') '  -2.12%  '
Set-TextValue $ws.Range('D49') '389.85'
Set-TextValue $ws.Range('E49') '  -0.12%  '
Set-TextValue $ws.Range('E50') '  -4.93%  '
Set-TextValue $ws.Range('D51') '25.74'
Set-TextValue $ws.Range('E51') '  -1.50%  '
